$d = $word.ActiveDocument

# 1. Append two trailing spaces to the first paragraph's text, and split
#    off a brand-new blank paragraph right after it (using ^p in the Find
#    replacement keeps the new paragraph free of any run, matching a
#    genuine Word "press Enter" split).
$find = $d.Content.Find
[void]$find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ^p", 2)

# 2. Append the red "(This is a change ... )" note to the end of the first
#    paragraph, as three separate runs (matching the source structure).
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$insertPos = $r.End - 1

$run1 = $d.Range($insertPos, $insertPos)
$run1.InsertAfter("(This is a change – Version for branch ")
$run1.Font.Color = 192

$run2 = $d.Range($run1.End, $run1.End)
$run2.InsertAfter("main")
$run2.Font.Color = 192

$run3 = $d.Range($run2.End, $run2.End)
$run3.InsertAfter(")")
$run3.Font.Color = 192
